# Update the "Förändrad" (Changed) date column (C) for rows 2-28:
# the stored serial date value moves from 45457 (2024-06-14) to 45458 (2024-06-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45457) {
        $cell.Value2 = 45458
    }
}
